$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Column width tweaks (col C and col J) ---
$ws.Columns.Item(3).ColumnWidth = 22.45
$ws.Columns.Item(10).ColumnWidth = 26.7

# --- Freeze panes: freeze column A and rows 1-3, active cell B4 ---
[void]$ws.Range("B4").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Row 14 : Frank Parks entry ---
$ws.Range("A14").Value = "Frank Parks"
$ws.Range("B14").Value = "8`" f/5 Newtonian"
$ws.Range("C14").Value = "Astrophysics Mach 1"
$ws.Range("D14").Value = "LiISA"
$ws.Range("E14").Value = "ATIK 460ex"
$ws.Range("F14").Value = "ATIK 414"
$ws.Range("G14").Value = "ACP with APCC"
$ws.Range("H14").Value = "Maxim DL"
$ws.Range("I14").Value = "Yes"
$ws.Range("J14").Value = "AstroPlanner + SkyX"
$ws.Range("K14").Value = "ISIS"

# --- Row 15 : Forrest Simms (Woody) entry ---
$ws.Range("A15").Value = "Forrest Simms (Woody)"
$ws.Range("B15").Value = "PlaneWave CDK14 + TV85"
$ws.Range("C15").Value = "Astro-Physics AP1100GTO"

$ws.Range("D15").Value = "LISA + Starlight Xpress`nFilter Wheel on TV85"
$ws.Range("D15").WrapText = $true
$ws.Range("D15").VerticalAlignment = -4160

$ws.Range("E15").Value = "ATIK414ex on LISA`nATIK 460ex on TV85`n"
$ws.Range("E15").WrapText = $true
$ws.Range("E15").VerticalAlignment = -4160

$ws.Range("F15").Value = "Lodestar x2"
$ws.Range("G15").Value = "MaximDL 6.20  + APCC Pro"
$ws.Range("H15").Value = "Maxim DL 6.20"
$ws.Range("I15").Value = "Yes"

$ws.Range("J15").Value = "AstroPlanner + Paolo Berardi`n Miles_Search excel spreadsheet"
$ws.Range("J15").WrapText = $true
$ws.Range("J15").VerticalAlignment = -4160

$ws.Range("K15").Value = "ISIS`nDemetra`nRspec"
$ws.Range("K15").WrapText = $true
$ws.Range("K15").VerticalAlignment = -4160

$ws.Rows.Item(15).RowHeight = 45
